# Journal-de-Travail_SchertenleibRomain.xlsx edit
# "remplire le JDT et faire le bilan"
#
# Fills in the "Journal de travail" sheet for the weeks of 05.05.2025
# (week 19), 12.05.2025 (week 20) and 19.05.2025 (week 21), and updates
# the window/view state to reflect where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# --- Week of 05.05.2025 (ISO week 19): dates were missing on rows 21-24,
#     the rest of the row (hours / type / description) was already filled. ---
$ws.Range("B21").Value = 45782
$ws.Range("B22").Value = 45782
$ws.Range("B23").Value = 45782
$ws.Range("B24").Value = 45782

# --- Week of 12.05.2025 (ISO week 20): rows 26-30, brand new entries. ---
$ws.Range("B26").Value = 45789
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = "Analyse"
$ws.Range("F26").Value = "reprise du projet "

$ws.Range("B27").Value = 45789
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 30
$ws.Range("E27").Value = "Implémentation"
$ws.Range("F29").Value = "documentation du le cache redis"
$ws.Range("F27").Value = "faire la partie 2 du projet "

$ws.Range("B28").Value = 45789
$ws.Range("D28").Value = 35
$ws.Range("E28").Value = "Implémentation"
$ws.Range("F28").Value = "finire l'implementation du todoscontrolleur "

$ws.Range("B29").Value = 45789
$ws.Range("D29").Value = 25
$ws.Range("E29").Value = "Documentation"

$ws.Range("B30").Value = 45789
$ws.Range("D30").Value = 15
$ws.Range("E30").Value = "Analyse"
$ws.Range("F30").Value = "faire des recherches sur l'implementation du todoscontolleur "

# --- Week of 19.05.2025 (ISO week 21): rows 32-37, brand new entries. ---
$ws.Range("B32").Value = 45796
$ws.Range("D32").Value = 10
$ws.Range("E32").Value = "Analyse"
$ws.Range("F32").Value = "reprise du projet "

$ws.Range("B33").Value = 45796
$ws.Range("C33").Value = 1
$ws.Range("E33").Value = "Implémentation"
$ws.Range("F33").Value = "corriger les erreurs de delete sur la partie 1"

$ws.Range("B34").Value = 45796
$ws.Range("D34").Value = 30
$ws.Range("E34").Value = "Implémentation"
$ws.Range("F34").Value = "corriger les erreurs d'update sur la partie 1"

$ws.Range("B35").Value = 45796
$ws.Range("D35").Value = 35
$ws.Range("E35").Value = "Implémentation"
$ws.Range("F35").Value = "corriger un point sur la partie 2 avec les utilisatueurs"

$ws.Range("B36").Value = 45796
$ws.Range("D36").Value = 45
$ws.Range("E36").Value = "Documentation"
$ws.Range("F36").Value = "documentation du le cache redis et tester des bouts de code "

$ws.Range("B37").Value = 45796
$ws.Range("D37").Value = 30
$ws.Range("E37").Value = "Implémentation"
$ws.Range("F37").Value = "commencer l'implémentation du cache redis"

# --- Window / view state: re-freeze the header rows (top 6 rows) and
#     scroll the frozen pane down so row 22 onward is visible, finishing
#     with the cursor on B37 (the last entry typed). ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A7").Select()
$win.FreezePanes = $true
$ws.Range("B37").Select()

